$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'26.237.09"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.Value = "'1.645.49"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "'  +0.51%  "
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.Value = "'  -0.14%  "
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.Value = "'216.86"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = "'  +0.64%  "
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.Value = "'  +0.46%  "
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.Value = "'  -0.18%  "
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.Value = "'0.258"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = "'  -0.29%  "
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.Value = "'  -0.28%  "
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.Value = "'19.99"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "'  +1.19%  "
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.Value = "'0.0793"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "'  +0.20%  "
$cell.Style = "Normal"
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$cell = $ws.Range("D12")
$cell.Value = "'4.30"
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.Value = "'  +0.36%  "
$cell.Style = "Normal"
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$cell = $ws.Range("D13")
$cell.Value = "'1.873.22"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = "'  +0.57%  "
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = "'1.629.48"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "'  -0.48%  "
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "'  -2.00%  "
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = "'  -0.60%  "
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.Value = "'63.56"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "'  +0.62%  "
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.Value = "'26.225.87"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "'  +1.56%  "
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "'  -0.18%  "
$cell.Style = "Normal"
$cell = $ws.Range("D20")
$cell.Value = "'195.72"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "'  +1.33%  "
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.Value = "'  -0.74%  "
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.Value = "'  +0.59%  "
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.Value = "'6.36"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = "'  -0.54%  "
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.Value = "'143.35"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = "'  +0.57%  "
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "'  -0.18%  "
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.Value = "'  -2.21%  "
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "'  +1.66%  "
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.Value = "'6.93"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.Value = "'  -0.27%  "
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.Value = "'15.61"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = "'  +0.39%  "
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.Value = "'  +1.25%  "
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.Value = "'0.0505"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.Value = "'  +2.21%  "
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.Value = "'  +0.37%  "
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = "'  +0.17%  "
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.Value = "'  +1.36%  "
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = "'  +1.30%  "
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.Value = "'  +1.01%  "
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.Value = "'1.135.84"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.Value = "'  +0.18%  "
$cell.Style = "Normal"
$cell = $ws.Range("D38")
$cell.Value = "'0.554"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = "'  +1.68%  "
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.Value = "'2.49"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.Value = "'  -1.51%  "
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.Value = "'0.0157"
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.Value = "'  +0.72%  "
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.Value = "'  -0.14%  "
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = "'  +1.83%  "
$cell.Style = "Normal"
$cell = $ws.Range("D43")
$cell.Value = "'100.16"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "'  -0.26%  "
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "'  -1.03%  "
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.Value = "'1.782.82"
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.Value = "'  +0.63%  "
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.Value = "'56.27"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "'  +1.75%  "
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.Value = "'1.49"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.Value = "'  +5.18%  "
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.Value = "'0.0518"
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.Value = "'  +2.87%  "
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = "'  +0.24%  "
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.Value = "'7.67"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.Value = "'  +2.37%  "
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.Value = "'0.0974"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.Value = "'  +2.20%  "
$cell.Style = "Normal"
